$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1021123333333333
$ws.Range("H2").Value = 0.306337
$ws.Range("I2").Value = 0.3579255212252356
$ws.Range("J2").Value = 0.3579255212252356
$ws.Range("M2").Value = 48.42420966666666
$ws.Range("N2").Value = 145.272629
$ws.Range("O2").Value = 0.6311762527593259
$ws.Range("P2").Value = 0.6311762527593258
$ws.Range("Q2").Value = 4.944709038885888
$ws.Range("R2").Value = 44.50238134997299
$ws.Range("S2").Value = 0.2259140892538728
$ws.Range("T2").Value = 0.2259140892538727
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1021123333333333
$ws.Range("H3").Value = 0.306337
$ws.Range("I3").Value = 0.3579255212252356
$ws.Range("J3").Value = 0.3579255212252356
$ws.Range("M3").Value = 6.849914666666667
$ws.Range("O3").Value = 0.08928392431779728
$ws.Range("P3").Value = 0.08928392431779726
$ws.Range("Q3").Value = 0.6994607697475554
$ws.Range("R3").Value = 6.295146927727999
$ws.Range("S3").Value = 0.03195699514848208
$ws.Range("T3").Value = 0.03195699514848208
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1021123333333333
$ws.Range("H4").Value = 0.306337
$ws.Range("I4").Value = 0.3579255212252356
$ws.Range("J4").Value = 0.3579255212252356
$ws.Range("N4").Value = 64.33937399999999
$ws.Range("O4").Value = 0.2795398229228769
$ws.Range("P4").Value = 0.2795398229228769
$ws.Range("Q4").Value = 2.189947868115333
$ws.Range("R4").Value = 19.709530813038
$ws.Range("S4").Value = 0.1000544368228808
$ws.Range("T4").Value = 0.1000544368228808
$ws.Range("I5").Value = 0.2397636084069039
$ws.Range("J5").Value = 0.2397636084069039
$ws.Range("M5").Value = 48.42420966666666
$ws.Range("N5").Value = 145.272629
$ws.Range("O5").Value = 0.6311762527593259
$ws.Range("P5").Value = 0.6311762527593258
$ws.Range("Q5").Value = 3.312312789619333
$ws.Range("R5").Value = 29.810815106574
$ws.Range("S5").Value = 0.151333095902324
$ws.Range("T5").Value = 0.151333095902324
$ws.Range("I6").Value = 0.2397636084069039
$ws.Range("J6").Value = 0.2397636084069039
$ws.Range("M6").Value = 6.849914666666667
$ws.Range("O6").Value = 0.08928392431779728
$ws.Range("P6").Value = 0.08928392431779726
$ws.Range("Q6").Value = 0.4685478630293334
$ws.Range("S6").Value = 0.02140703586716399
$ws.Range("T6").Value = 0.02140703586716398
$ws.Range("I7").Value = 0.2397636084069039
$ws.Range("J7").Value = 0.2397636084069039
$ws.Range("N7").Value = 64.33937399999999
$ws.Range("O7").Value = 0.2795398229228769
$ws.Range("P7").Value = 0.2795398229228769
$ws.Range("S7").Value = 0.06702347663741591
$ws.Range("T7").Value = 0.06702347663741591
$ws.Range("I8").Value = 0.4023108703678605
$ws.Range("J8").Value = 0.4023108703678605
$ws.Range("M8").Value = 48.42420966666666
$ws.Range("N8").Value = 145.272629
$ws.Range("O8").Value = 0.6311762527593259
$ws.Range("P8").Value = 0.6311762527593258
$ws.Range("Q8").Value = 5.557888664491666
$ws.Range("R8").Value = 50.02099798042499
$ws.Range("S8").Value = 0.2539290676031291
$ws.Range("T8").Value = 0.2539290676031291
$ws.Range("I9").Value = 0.4023108703678605
$ws.Range("J9").Value = 0.4023108703678605
$ws.Range("M9").Value = 6.849914666666667
$ws.Range("O9").Value = 0.08928392431779728
$ws.Range("P9").Value = 0.08928392431779726
$ws.Range("Q9").Value = 0.7861989558666667
$ws.Range("R9").Value = 7.0757906028
$ws.Range("S9").Value = 0.0359198933021512
$ws.Range("T9").Value = 0.03591989330215119
$ws.Range("I10").Value = 0.4023108703678605
$ws.Range("J10").Value = 0.4023108703678605
$ws.Range("N10").Value = 64.33937399999999
$ws.Range("O10").Value = 0.2795398229228769
$ws.Range("P10").Value = 0.2795398229228769
$ws.Range("Q10").Value = 2.461517216949999
$ws.Range("S10").Value = 0.1124619094625802
$ws.Range("T10").Value = 0.1124619094625802
